$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row for Delaware into Region 3 (after Puerto Rico, before District of Columbia)
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "United States"
$ws.Range("B12").Value = "Region 3"
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = "Delaware"

# 2) Fix the renamed / typo'd state names
$ws.Range("D11").Value = "Virgin Islands"
$ws.Range("D36").Value = "Texas"

# 3) Recompute the HHSRegionNumber column (C) so it consistently matches the
#    region number parsed out of the HHSRegionName column (B) - this corrects
#    the off-by-one bug where the count incremented a row too early.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $regionName = $ws.Cells.Item($r, 2).Value2
    $regionNumber = [int]($regionName -replace '[^0-9]', '')
    $ws.Cells.Item($r, 3).Value = $regionNumber
}

# 4) Restore the selected cell as recorded after the edit
$ws.Range("D12").Select() | Out-Null
